# Cronograma del proyecto: corre las fechas de inicio (columna J, filas 18-26)
# tres dias (de 45407 a 45410) y ajusta la duracion de la tarea de la fila 38
# (de I38+8 a I38+4). El resto de fechas de la hoja son formulas que dependen
# de estas celdas y se recalculan automaticamente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 18..26) {
    $ws.Cells.Item($row, 10).Value = 45410
}

$ws.Range("J38").Formula = "=I38+4"

$wb.Application.CalculateFull()
